# Weekly update: insert two new daily-price rows for "Vega Monumental
# Concepción - Cebolla" (matches commit message "Fruta / hortaliza, semanal").
#
# The new records are inserted immediately above what is currently row 261,
# pushing the existing rows 261-290 down to 263-292 (all their content,
# formatting and the sheet's used range grow accordingly).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows (A:R) above row 261; existing data shifts down.
$ws.Range("A261:R262").Insert()

# --- New row 261: "1a nueva(o)" ---
$ws.Range("A261").Value = 11
$ws.Range("B261").Value = "Vega Monumental Concepción"
$ws.Range("C261").Value = "Bíobío"
$ws.Range("D261").Value = 44505
$ws.Range("E261").Value = 8
$ws.Range("F261").Value = 100112004
$ws.Range("G261").Value = "Cebolla"
$ws.Range("H261").Value = "Sin especificar"
$ws.Range("I261").Value = "1a nueva(o)"
$ws.Range("J261").Value = 430
$ws.Range("K261").Value = 5000
$ws.Range("L261").Value = 5500
$ws.Range("M261").Value = 5291
$ws.Range("N261").Value = "`$/malla 18 kilos"
$ws.Range("O261").Value = "Región de O'Higgins"
$ws.Range("P261").Value = 294
$ws.Range("Q261").Value = 18
$ws.Range("R261").Value = "Hortaliza"

# --- New row 262: "2a nueva(o)" ---
$ws.Range("A262").Value = 11
$ws.Range("B262").Value = "Vega Monumental Concepción"
$ws.Range("C262").Value = "Bíobío"
$ws.Range("D262").Value = 44505
$ws.Range("E262").Value = 8
$ws.Range("F262").Value = 100112004
$ws.Range("G262").Value = "Cebolla"
$ws.Range("H262").Value = "Sin especificar"
$ws.Range("I262").Value = "2a nueva(o)"
$ws.Range("J262").Value = 200
$ws.Range("K262").Value = 4500
$ws.Range("L262").Value = 4500
$ws.Range("M262").Value = 4500
$ws.Range("N262").Value = "`$/malla 18 kilos"
$ws.Range("O262").Value = "Región de O'Higgins"
$ws.Range("P262").Value = 250
$ws.Range("Q262").Value = 18
$ws.Range("R262").Value = "Hortaliza"
